# Update cryptocurrency price/volume data (and a few name/link swaps)
# to reflect the latest scrape, per commit:
# "Updated cryptos list on Thu Dec 28 06:31:16 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.184.69"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "2.380.95"
$ws.Range("E3").Value = "  +6.93%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "321.86"
$ws.Range("E5").Value = "  +9.64%  "
$ws.Range("D6").Value = "106.51"
$ws.Range("E6").Value = "  -5.52%  "
$ws.Range("D7").Value = "0.648"
$ws.Range("E7").Value = "  +4.61%  "
$ws.Range("D9").Value = "0.646"
$ws.Range("E9").Value = "  +7.47%  "
$ws.Range("D10").Value = "42.16"
$ws.Range("E10").Value = "  -3.53%  "
$ws.Range("D11").Value = "0.0938"
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").Value = "17.27"
$ws.Range("E13").Value = "  +15.35%  "
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "2.751.17"
$ws.Range("E16").Value = "  +7.41%  "
$ws.Range("D17").Value = "2.388.32"
$ws.Range("E17").Value = "  +6.85%  "
$ws.Range("D18").Value = "43.217.96"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").Value = "7.48"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").Value = "0.0000108"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").Value = "76.17"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").Value = "271.07"
$ws.Range("E22").Value = "  +15.04%  "
$ws.Range("D23").Value = "3.41"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "9.84"
$ws.Range("E25").Value = "  +10.49%  "
$ws.Range("D26").Value = "11.79"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "22.83"
$ws.Range("E28").Value = "  +7.06%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("D29").Value = "38.17"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("B30").Value = "Monero"
$ws.Range("D30").Value = "176.88"
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "3.21"
$ws.Range("D33").Value = "0.0927"
$ws.Range("E33").Value = "  +4.86%  "
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("E35").Value = "  +5.75%  "
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").Value = "4.16"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("D38").Value = "0.0367"
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").Value = "0.107"
$ws.Range("E39").Value = "  +3.22%  "
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +17.18%  "
$ws.Range("D41").Value = "1.58"
$ws.Range("E41").Value = "  +20.51%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("D42").Value = "0.232"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("B43").Value = "Aave"
$ws.Range("D43").Value = "123.87"
$ws.Range("E43").Value = "  +22.51%  "
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "69.17"
$ws.Range("E44").Value = "  -3.70%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "12.52"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("D47").Value = "91.87"
$ws.Range("E47").Value = "  +67.59%  "
$ws.Range("D48").Value = "9.52"
$ws.Range("E48").Value = "  +12.91%  "
$ws.Range("E49").Value = "  +4.87%  "
$ws.Range("D50").Value = "1.31"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  +5.59%  "
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"